# PowerShell-style Excel COM-interop script
# Applies update to Notified_Production_Wind.xlsx:
#  - Column A (timestamps, rows 2-97): shift each date serial by +11 days
#  - Column B (notified production values, rows 2-97): replace with new values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aVals = @(45797.01041666666,45797.02083333334,45797.03125,45797.04166666666,45797.05208333334,45797.0625,45797.07291666666,45797.08333333334,45797.09375,45797.10416666666,45797.11458333334,45797.125,45797.13541666666,45797.14583333334,45797.15625,45797.16666666666,45797.17708333334,45797.1875,45797.19791666666,45797.20833333334,45797.21875,45797.22916666666,45797.23958333334,45797.25,45797.26041666666,45797.27083333334,45797.28125,45797.29166666666,45797.30208333334,45797.3125,45797.32291666666,45797.33333333334,45797.34375,45797.35416666666,45797.36458333334,45797.375,45797.38541666666,45797.39583333334,45797.40625,45797.41666666666,45797.42708333334,45797.4375,45797.44791666666,45797.45833333334,45797.46875,45797.47916666666,45797.48958333334,45797.5,45797.51041666666,45797.52083333334,45797.53125,45797.54166666666,45797.55208333334,45797.5625,45797.57291666666,45797.58333333334,45797.59375,45797.60416666666,45797.61458333334,45797.625,45797.63541666666,45797.64583333334,45797.65625,45797.66666666666,45797.67708333334,45797.6875,45797.69791666666,45797.70833333334,45797.71875,45797.72916666666,45797.73958333334,45797.75,45797.76041666666,45797.77083333334,45797.78125,45797.79166666666,45797.80208333334,45797.8125,45797.82291666666,45797.83333333334,45797.84375,45797.85416666666,45797.86458333334,45797.875,45797.88541666666,45797.89583333334,45797.90625,45797.91666666666,45797.92708333334,45797.9375,45797.94791666666,45797.95833333334,45797.96875,45797.97916666666,45797.98958333334,45798)
$bVals = @(867,867,872,870,815,808,804,799,810,811,818,820,837,837,834,839,829,830,831,824,834,832,829,829,832,829,828,817,848,842,849,859,1103,1111,1120,1126,1393,1392,1391,1390,1307,1304,1301,1298,1181,1179,1177,1174,1106,1105,1104,1103,1102,1101,1101,1100,1080,1079,1078,1078,895,894,893,892,805,803,802,800,722,717,712,708,579,575,572,569,459,456,453,451,360,358,357,356,282,279,280,280,239,238,238,238,0,0,0,0)

$n = $aVals.Length
$startRow = 2

for ($i = 0; $i -lt $n; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}
